# Update F-column (想去人数 / want-to-go count) values across all sheets
# to match the refreshed data snapshot (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 848
$ws1.Range("F3").Value = 1750
$ws1.Range("F7").Value = 1365
$ws1.Range("F8").Value = 2065
$ws1.Range("F9").Value = 965
$ws1.Range("F14").Value = 3907
$ws1.Range("F16").Value = 365
$ws1.Range("F17").Value = 2982
$ws1.Range("F18").Value = 791
$ws1.Range("F21").Value = 108
$ws1.Range("F22").Value = 2044
$ws1.Range("F23").Value = 1171
$ws1.Range("F24").Value = 1850
$ws1.Range("F25").Value = 382
$ws1.Range("F26").Value = 204
$ws1.Range("F27").Value = 10
$ws1.Range("F28").Value = 8322
$ws1.Range("F29").Value = 5627
$ws1.Range("F30").Value = 352
$ws1.Range("F31").Value = 172
$ws1.Range("F32").Value = 748
$ws1.Range("F33").Value = 762
$ws1.Range("F34").Value = 3476
$ws1.Range("F37").Value = 381
$ws1.Range("F38").Value = 33
$ws1.Range("F39").Value = 188
$ws1.Range("F40").Value = 154
$ws1.Range("F41").Value = 4613
$ws1.Range("F42").Value = 841
$ws1.Range("F43").Value = 71
$ws1.Range("F44").Value = 392

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F15").Value = 106

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 8246
$ws3.Range("F3").Value = 358
$ws3.Range("F4").Value = 1264

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 848
$ws4.Range("F3").Value = 358
$ws4.Range("F4").Value = 1264
$ws4.Range("F6").Value = 1750
$ws4.Range("F9").Value = 1365
$ws4.Range("F10").Value = 2065
$ws4.Range("F11").Value = 965
$ws4.Range("F15").Value = 3907
$ws4.Range("F16").Value = 365
$ws4.Range("F17").Value = 2982
$ws4.Range("F18").Value = 791
$ws4.Range("F21").Value = 2044
$ws4.Range("F27").Value = 1850
$ws4.Range("F28").Value = 106
$ws4.Range("F29").Value = 204
$ws4.Range("F30").Value = 10
$ws4.Range("F31").Value = 8322
$ws4.Range("F32").Value = 5627
$ws4.Range("F34").Value = 352
$ws4.Range("F35").Value = 172
$ws4.Range("F36").Value = 748
$ws4.Range("F37").Value = 762
$ws4.Range("F40").Value = 381
$ws4.Range("F42").Value = 154
$ws4.Range("F43").Value = 4613
$ws4.Range("F44").Value = 841
$ws4.Range("F45").Value = 392
